# Fix a language error ("taalfout") in the "Tekst alcohol" placeholder that
# lives on the slide layout used by slide 7 ("Alcohol - 2"):
#   "... ten opzichte van hangt samen ..."
# should read
#   "... ten opzichte van alcohol hangt samen ..."
#
# The text lives on the slide's CustomLayout (slide master layout), not on
# the slide itself, so we have to walk each slide's CustomLayout shapes to
# find and patch it in place.

$p = $ppt.ActivePresentation

$newText = "De houding van ouders ten opzichte van alcohol hangt samen met het alcoholgebruik van jongeren. Van de jongeren die geen (of alleen een slokje) alcohol drinken geeft xx% aan dat hun ouders het goed zouden vinden als ze toch alcohol zouden drinken. Van de jongeren die alcohol drinken geeft xx% aan dat hun ouders alcohol drinken goed vinden."
$needle = "ten opzichte van hangt samen met het alcoholgebruik"

$patched = 0
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $cl = $s.CustomLayout
    for ($k = 1; $k -le $cl.Shapes.Count; $k++) {
        $sh = $cl.Shapes.Item($k)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text.Contains($needle)) {
                $tr.Text = $newText
                $patched = $patched + 1
            }
        }
    }
}
Write-Host "Patched shapes:" $patched
